$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the burn-down values in column C: the remaining work dropped
# again on 12/5 (row 12) and again on 12/8 (row 14). The rest of the
# shared "=prev" chain then recalculates to the new totals automatically.
$ws.Range("C12").Formula = "=C11-10"
$ws.Range("C14").Formula = "=C13-6"

# Move the active selection to C14, matching where the edit was made.
$ws.Range("C14").Select()
